# Insert a new weekly price record at row 232 of the single data sheet.
#
# The sheet holds one price observation per row (rows 2..349, with a header
# in row 1). This edit inserts one additional observation right before the
# existing row 232, which pushes every subsequent row down by one (old row
# 232 -> new row 233, ..., old row 349 -> new row 350). The new row reuses
# the price-range / unit / origin metadata (columns A,B,C,E,F,G,H,I,K,L,M,
# N,O,P,Q,R) of the row it displaces, only the date (D) and volume (J)
# columns get genuinely new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 232
$lastCol = 18   # columns A..R

# Shift row $newRow (and everything below it) down by one row, leaving a
# blank row $newRow behind with the same row-level formatting Excel assigns
# on insert (e.g. the date style that was on column D).
$ws.Rows.Item($newRow).Insert()

# The data that used to live at $newRow now lives at $newRow + 1. Copy its
# values back up into the freshly inserted row so the new record starts
# out identical to the one it is based on.
for ($col = 1; $col -le $lastCol; $col++) {
    $ws.Cells.Item($newRow, $col).Value = $ws.Cells.Item($newRow + 1, $col).Value2()
}

# Now overwrite just the two columns that actually differ for the new
# observation: Fecha (D, column 4) and Volumen (J, column 10).
$ws.Cells.Item($newRow, 4).Value = 44917
$ws.Cells.Item($newRow, 10).Value = 1360
